# Generate Report for Archive
#
# The localization status report is refreshed for archival:
#  - the in-flight items' Status moves from "Ready for handoff" to
#    "In Translation" on every sheet that surfaces it (Overview's per-locale
#    status columns, and each locale sheet's own Status column)
#  - the now-narrower Status column text lets the Status column(s) be
#    narrowed accordingly on every sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn / de-de sheets: column C is the "Status" column
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status column(s) to match the new, shorter text ---
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
